# Rename speaker "Davis" to "T" in the Speaker column (column D) of the DataSheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)  # Column D = Speaker
    if ($cell.Value2 -eq "Davis") {
        $cell.Value2 = "T"
    }
}
